$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected Diebold-Mariano statistics (DM_Stat), p-values (P_Value)
# and recomputed significance flag (Significativo) for each comparison row.

# Row 2: N_Calib_1=20, N_Calib_2=40
$ws.Range("C2").Value = -1.759873895633086
$ws.Range("D2").Value = 0.08742475010839024
$ws.Range("G2").Value = "No"

# Row 3: N_Calib_1=20, N_Calib_2=60
$ws.Range("C3").Value = -0.160181847350048
$ws.Range("D3").Value = 0.8736857220345882
$ws.Range("G3").Value = "No"

# Row 4: N_Calib_1=20, N_Calib_2=100
$ws.Range("C4").Value = 0.5000231264769484
$ws.Range("D4").Value = 0.6202797860514058
$ws.Range("G4").Value = "No"

# Row 5: N_Calib_1=20, N_Calib_2=200
$ws.Range("C5").Value = -0.8489520434153063
$ws.Range("D5").Value = 0.401845933985733
$ws.Range("G5").Value = "No"

# Row 6: N_Calib_1=40, N_Calib_2=60
$ws.Range("C6").Value = 1.793963106982317
$ws.Range("D6").Value = 0.08171671630001009
$ws.Range("G6").Value = "No"

# Row 7: N_Calib_1=40, N_Calib_2=100
$ws.Range("C7").Value = 2.182839810154933
$ws.Range("D7").Value = 0.03604886797550466
$ws.Range("G7").Value = "Sí"

# Row 8: N_Calib_1=40, N_Calib_2=200
$ws.Range("C8").Value = 1.596281245383899
$ws.Range("D8").Value = 0.1196802553859901
$ws.Range("G8").Value = "No"

# Row 9: N_Calib_1=60, N_Calib_2=100
$ws.Range("C9").Value = 0.7514717581631893
$ws.Range("D9").Value = 0.4575392062173886
$ws.Range("G9").Value = "No"

# Row 10: N_Calib_1=60, N_Calib_2=200
$ws.Range("C10").Value = -0.4553226048435484
$ws.Range("D10").Value = 0.6517716420648192
$ws.Range("G10").Value = "No"

# Row 11: N_Calib_1=100, N_Calib_2=200
$ws.Range("C11").Value = -1.012989453846747
$ws.Range("D11").Value = 0.3182201180564375
$ws.Range("G11").Value = "No"
